$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Julio de 2020 a las 10:40"

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 1389170
$ws.Range("C6").Value = 3676
$ws.Range("D6").Value = 887296
$ws.Range("E6").Value = 469746
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 32
$ws.Range("H6").Value = 32128

# Row 7: Rusia
$ws.Range("A7").Value = "Rusia"
$ws.Range("B7").Value = 812485
$ws.Range("C7").Value = 5765
$ws.Range("D7").Value = 600250
$ws.Range("E7").Value = 198966
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 77
$ws.Range("H7").Value = 13269

# Row 31: Filipinas
$ws.Range("A31").Value = "Filipinas"
$ws.Range("B31").Value = 80448
$ws.Range("C31").Value = 2036
$ws.Range("D31").Value = 26110
$ws.Range("E31").Value = 52406
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 35
$ws.Range("H31").Value = 1932

# Row 32: Ecuador
$ws.Range("A32").Value = "Ecuador"
$ws.Range("B32").Value = 80036
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 34544
$ws.Range("E32").Value = 39985
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 5507

# Row 33: Suecia
$ws.Range("A33").Value = "Suecia"
$ws.Range("B33").Value = 78997
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 5697

# Row 38: Ucrania
$ws.Range("A38").Value = "Ucrania"
$ws.Range("B38").Value = 64849
$ws.Range("C38").Value = 920
$ws.Range("D38").Value = 35807
$ws.Range("E38").Value = 27437
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 15
$ws.Range("H38").Value = 1605

# Row 40: Israel
$ws.Range("A40").Value = "Israel"
$ws.Range("B40").Value = 61388
$ws.Range("C40").Value = 710
$ws.Range("D40").Value = 26959
$ws.Range("E40").Value = 33965
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 7
$ws.Range("H40").Value = 464

# Row 41: Republica Dominicana
$ws.Range("A41").Value = "Republica Dominicana"
$ws.Range("B41").Value = 60896
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 27980
$ws.Range("E41").Value = 31861
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 1055

# Row 45: Singapur
$ws.Range("A45").Value = "Singapur"
$ws.Range("B45").Value = 50369
$ws.Range("C45").Value = 481
$ws.Range("D45").Value = 45352
$ws.Range("E45").Value = 4990
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 27

# Row 46: Portugal
$ws.Range("A46").Value = "Portugal"
$ws.Range("B46").Value = 49955
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 35010
$ws.Range("E46").Value = 13229
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 1716

# Row 51: Barein
$ws.Range("A51").Value = "Barein"
$ws.Range("B51").Value = 38747
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 35205
$ws.Range("E51").Value = 3403
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 139

# Row 53: Armenia
$ws.Range("A53").Value = "Armenia"
$ws.Range("B53").Value = 37317
$ws.Range("C53").Value = 321
$ws.Range("D53").Value = 26478
$ws.Range("E53").Value = 10134
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 5
$ws.Range("H53").Value = 705

# Row 63: Moldavia
$ws.Range("A63").Value = "Moldavia"
$ws.Range("B63").Value = 22828
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 15909
$ws.Range("E63").Value = 6186
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 733

# Row 64: Austria
$ws.Range("A64").Value = "Austria"
$ws.Range("B64").Value = 20472
$ws.Range("C64").Value = 134
$ws.Range("D64").Value = 18209
$ws.Range("E64").Value = 1551
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 712

# Row 123: Eslovaquia
$ws.Range("A123").Value = "Eslovaquia"
$ws.Range("B123").Value = 2179
$ws.Range("C123").Value = 38
$ws.Range("D123").Value = 1577
$ws.Range("E123").Value = 574
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 28

# Row 124: Suazilandia
$ws.Range("A124").Value = "Suazilandia"
$ws.Range("B124").Value = 2142
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 940
$ws.Range("E124").Value = 1174
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 28

# Row 127: Lituania
$ws.Range("A127").Value = "Lituania"
$ws.Range("B127").Value = 2008
$ws.Range("C127").Value = 7
$ws.Range("D127").Value = 1616
$ws.Range("E127").Value = 312
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 80

# Row 139: Letonia
$ws.Range("A139").Value = "Letonia"
$ws.Range("B139").Value = 1219
$ws.Range("C139").Value = 13
$ws.Range("D139").Value = 1045
$ws.Range("E139").Value = 143
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 31

# Row 210: Islas Malvinas
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("B210").Value = 13
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 13
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Row 211: Groenlandia
$ws.Range("A211").Value = "Groenlandia"
$ws.Range("B211").Value = 13
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 13
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0
